$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.991.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.87%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.87%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6229'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.32%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07520'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.57%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2913'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.86%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07723'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.830.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.95%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.937'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6633'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.85%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001004'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +14.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.41'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.019'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.82%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.977.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '225.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.138'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9995'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.450'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1369'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.36%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.490'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.71%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.081'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.48%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.017'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.192'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05187'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.843'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7370'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.49%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.137'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.34%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.695'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.82%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.241.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.759'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.03%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01783'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.90%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.317'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.41%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8950'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.60%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9996'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.24%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.980.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.72%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.23%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5103'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.75%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.46%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4021'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.825'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05742'
$ws.Range('D50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.634'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.97%  '
